$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that look like numbers remain stored as text
# (preserves values such as "1.002", "316.02", "1.001" as strings, matching
# the source data which uses dotted/decimal notation as plain text).
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D8:D21").NumberFormat = "@"
$ws.Range("D23:D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.690.76"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "1.694.18"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "316.02"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "1.488"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "1.003"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "53.11"
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").Value = "0.08851"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "7.267"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "23.65"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "8.058"
$ws.Range("E15").Value = "  +8.67%  "
$ws.Range("D16").Value = "0.00001319"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "1.695.64"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "100.05"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "0.07010"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "19.56"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "7.047"
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "14.34"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "24.680.63"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "3.269"
$ws.Range("E25").Value = "  +10.27%  "
$ws.Range("D26").Value = "2.358"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D27").Value = "22.78"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("D28").Value = "163.61"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "136.49"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "5.180"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "7.494"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "1.879.09"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "1.072"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "0.08599"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "7.145"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").Value = "11.41"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "0.2750"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").Value = "1.926"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "14.46"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "0.09222"
$ws.Range("E40").Value = "  +3.12%  "
$ws.Range("D41").Value = "0.02727"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "1.464"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "0.7672"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("D44").Value = "16.02"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").Value = "0.7193"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "2.580"
$ws.Range("E46").Value = "  +5.99%  "
$ws.Range("D47").Value = "4.216"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "1.323"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "139.56"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "0.07991"
$ws.Range("E51").Value = "  +0.91%  "
